# wip for new file path code
#
# Adds a new "PrimaryTopLevelFolderType"-style flag column (G) to Sheet1,
# marking most existing rows (3-17 and 19) with "yes", then leaves the
# selection where the author's cursor ended up (G20, the first empty cell
# below the new column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G: flag every data row with "yes" except the header rows (1-2)
# and rows 18/20, matching the source edit exactly.
$ws.Range("G3:G17").Value = "yes"
$ws.Range("G19").Value = "yes"

# Column G is otherwise default width; column I picks up a wider, bestFit-like
# width as a side effect of the author's edits further along that row.
$ws.Columns.Item(9).ColumnWidth = 57.1640625

# Leave the cursor on G20 (first empty cell under the new column), same as
# the recorded selection in the saved file.
$ws.Range("G20").Select() | Out-Null
